# Update the "取得日時" (retrieved datetime) column for all existing data
# rows on the "ランサーズ" sheet to reflect the latest append timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-09-30 01:40:31"

# Data rows are 2 through 13 (row 1 is the header "取得日時").
for ($row = 2; $row -le 13; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
